# Apply updated PtX demand categories (Fossil Gases / Fossil Liquids added, rows re-sequenced)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Fossil Gases (2030)
$ws.Cells.Item(7, 1).Value = "Fossil Gases"
$ws.Cells.Item(7, 2).Value = 2030
$ws.Cells.Item(7, 3).Value = $null
$ws.Cells.Item(7, 4).Value = $null
$ws.Cells.Item(7, 5).Value = $null
$ws.Cells.Item(7, 6).Value = 0.001424263216952223
$ws.Cells.Item(7, 7).Value = $null
$ws.Cells.Item(7, 8).Value = $null
$ws.Cells.Item(7, 9).Value = 0.00008025862017923084
$ws.Cells.Item(7, 10).Value = $null
$ws.Cells.Item(7, 11).Value = $null

# Row 8: Synthetic Liquids (2030)
$ws.Cells.Item(8, 1).Value = "Synthetic Liquids"
$ws.Cells.Item(8, 2).Value = 2030
$ws.Cells.Item(8, 3).Value = $null
$ws.Cells.Item(8, 4).Value = $null
$ws.Cells.Item(8, 5).Value = $null
$ws.Cells.Item(8, 6).Value = $null
$ws.Cells.Item(8, 7).Value = $null
$ws.Cells.Item(8, 8).Value = $null
$ws.Cells.Item(8, 9).Value = $null
$ws.Cells.Item(8, 10).Value = $null
$ws.Cells.Item(8, 11).Value = $null

# Row 9: Biogenic Liquids (2030)
$ws.Cells.Item(9, 1).Value = "Biogenic Liquids"
$ws.Cells.Item(9, 2).Value = 2030
$ws.Cells.Item(9, 3).Value = $null
$ws.Cells.Item(9, 4).Value = $null
$ws.Cells.Item(9, 5).Value = $null
$ws.Cells.Item(9, 6).Value = 0.006658482628641779
$ws.Cells.Item(9, 7).Value = 0.00003393298861638786
$ws.Cells.Item(9, 8).Value = 0.00220440146393377
$ws.Cells.Item(9, 9).Value = 0.0037880844391939
$ws.Cells.Item(9, 10).Value = 0.00009263140517044296
$ws.Cells.Item(9, 11).Value = 0.0027623925437353

# Row 10: Fossil Liquids (2030)
$ws.Cells.Item(10, 1).Value = "Fossil Liquids"
$ws.Cells.Item(10, 2).Value = 2030
$ws.Cells.Item(10, 3).Value = $null
$ws.Cells.Item(10, 4).Value = $null
$ws.Cells.Item(10, 5).Value = $null
$ws.Cells.Item(10, 6).Value = 0.066564228137772
$ws.Cells.Item(10, 7).Value = 0.0002389800881419
$ws.Cells.Item(10, 8).Value = 0.0201243571619895
$ws.Cells.Item(10, 9).Value = 0.0239554169299955
$ws.Cells.Item(10, 10).Value = 0.0005606980903764
$ws.Cells.Item(10, 11).Value = 0.0255901401576523

# Row 11: Biomass [Solid] (2030)
$ws.Cells.Item(11, 1).Value = "Biomass [Solid]"
$ws.Cells.Item(11, 2).Value = 2030
$ws.Cells.Item(11, 3).Value = $null
$ws.Cells.Item(11, 4).Value = $null
$ws.Cells.Item(11, 5).Value = 0.001767050420365749
$ws.Cells.Item(11, 6).Value = $null
$ws.Cells.Item(11, 7).Value = $null
$ws.Cells.Item(11, 8).Value = $null
$ws.Cells.Item(11, 9).Value = $null
$ws.Cells.Item(11, 10).Value = $null
$ws.Cells.Item(11, 11).Value = $null

# Row 12: Renewable Energy Carrier (2030)
$ws.Cells.Item(12, 1).Value = "Renewable Energy Carrier"
$ws.Cells.Item(12, 2).Value = 2030
$ws.Cells.Item(12, 3).Value = $null
$ws.Cells.Item(12, 4).Value = $null
$ws.Cells.Item(12, 5).Value = 0.0006762847930484331
$ws.Cells.Item(12, 6).Value = $null
$ws.Cells.Item(12, 7).Value = $null
$ws.Cells.Item(12, 8).Value = $null
$ws.Cells.Item(12, 9).Value = $null
$ws.Cells.Item(12, 10).Value = $null
$ws.Cells.Item(12, 11).Value = $null

# Row 13: Overall Demand (2030)
$ws.Cells.Item(13, 1).Value = "Overall Demand"
$ws.Cells.Item(13, 2).Value = 2030
$ws.Cells.Item(13, 3).Value = $null
$ws.Cells.Item(13, 4).Value = 0.004519281559255572
$ws.Cells.Item(13, 5).Value = 0.003291958720677564
$ws.Cells.Item(13, 6).Value = 0.07512321926574525
$ws.Cells.Item(13, 7).Value = 0.0002729130767582878
$ws.Cells.Item(13, 8).Value = 0.02232875947830612
$ws.Cells.Item(13, 9).Value = 0.0279503680024534
$ws.Cells.Item(13, 10).Value = 0.0006533294955468429
$ws.Cells.Item(13, 11).Value = 0.0283525327013876

# Row 14: Hydrogen (2040)
$ws.Cells.Item(14, 1).Value = "Hydrogen"
$ws.Cells.Item(14, 2).Value = 2040
$ws.Cells.Item(14, 3).Value = $null
$ws.Cells.Item(14, 4).Value = $null
$ws.Cells.Item(14, 5).Value = $null
$ws.Cells.Item(14, 6).Value = 0.001781868501600853
$ws.Cells.Item(14, 7).Value = $null
$ws.Cells.Item(14, 8).Value = 0.00000007135385940852036
$ws.Cells.Item(14, 9).Value = 0.0001561632335787969
$ws.Cells.Item(14, 10).Value = $null
$ws.Cells.Item(14, 11).Value = $null

# Row 15: Methanol (2040)
$ws.Cells.Item(15, 1).Value = "Methanol"
$ws.Cells.Item(15, 2).Value = 2040
$ws.Cells.Item(15, 3).Value = $null
$ws.Cells.Item(15, 4).Value = $null
$ws.Cells.Item(15, 5).Value = $null
$ws.Cells.Item(15, 6).Value = $null
$ws.Cells.Item(15, 7).Value = $null
$ws.Cells.Item(15, 8).Value = $null
$ws.Cells.Item(15, 9).Value = $null
$ws.Cells.Item(15, 10).Value = $null
$ws.Cells.Item(15, 11).Value = $null

# Row 16: Ammonia (2040)
$ws.Cells.Item(16, 1).Value = "Ammonia"
$ws.Cells.Item(16, 2).Value = 2040
$ws.Cells.Item(16, 3).Value = $null
$ws.Cells.Item(16, 4).Value = 0.004453516139235275
$ws.Cells.Item(16, 5).Value = $null
$ws.Cells.Item(16, 6).Value = $null
$ws.Cells.Item(16, 7).Value = $null
$ws.Cells.Item(16, 8).Value = $null
$ws.Cells.Item(16, 9).Value = $null
$ws.Cells.Item(16, 10).Value = $null
$ws.Cells.Item(16, 11).Value = $null

# Row 17: Synthetic Gases (2040)
$ws.Cells.Item(17, 1).Value = "Synthetic Gases"
$ws.Cells.Item(17, 2).Value = 2040
$ws.Cells.Item(17, 3).Value = $null
$ws.Cells.Item(17, 4).Value = $null
$ws.Cells.Item(17, 5).Value = $null
$ws.Cells.Item(17, 6).Value = 0.0000000006452770915871963
$ws.Cells.Item(17, 7).Value = $null
$ws.Cells.Item(17, 8).Value = $null
$ws.Cells.Item(17, 9).Value = 0.00000000004887106473705467
$ws.Cells.Item(17, 10).Value = $null
$ws.Cells.Item(17, 11).Value = $null

# Row 18: Biogenic Gases (2040)
$ws.Cells.Item(18, 1).Value = "Biogenic Gases"
$ws.Cells.Item(18, 2).Value = 2040
$ws.Cells.Item(18, 3).Value = $null
$ws.Cells.Item(18, 4).Value = $null
$ws.Cells.Item(18, 5).Value = 0.003402768728558265
$ws.Cells.Item(18, 6).Value = 0.0001325854197537344
$ws.Cells.Item(18, 7).Value = $null
$ws.Cells.Item(18, 8).Value = $null
$ws.Cells.Item(18, 9).Value = 0.00003364620598855328
$ws.Cells.Item(18, 10).Value = $null
$ws.Cells.Item(18, 11).Value = $null

# Row 19: Fossil Gases (2040)
$ws.Cells.Item(19, 1).Value = "Fossil Gases"
$ws.Cells.Item(19, 2).Value = 2040
$ws.Cells.Item(19, 3).Value = $null
$ws.Cells.Item(19, 4).Value = $null
$ws.Cells.Item(19, 5).Value = $null
$ws.Cells.Item(19, 6).Value = 0.0007658853692322342
$ws.Cells.Item(19, 7).Value = $null
$ws.Cells.Item(19, 8).Value = $null
$ws.Cells.Item(19, 9).Value = 0.00008525818792382284
$ws.Cells.Item(19, 10).Value = $null
$ws.Cells.Item(19, 11).Value = $null

# Row 20: Synthetic Liquids (2040)
$ws.Cells.Item(20, 1).Value = "Synthetic Liquids"
$ws.Cells.Item(20, 2).Value = 2040
$ws.Cells.Item(20, 3).Value = $null
$ws.Cells.Item(20, 4).Value = $null
$ws.Cells.Item(20, 5).Value = $null
$ws.Cells.Item(20, 6).Value = $null
$ws.Cells.Item(20, 7).Value = $null
$ws.Cells.Item(20, 8).Value = $null
$ws.Cells.Item(20, 9).Value = $null
$ws.Cells.Item(20, 10).Value = $null
$ws.Cells.Item(20, 11).Value = $null

# Row 21: Biogenic Liquids (2040)
$ws.Cells.Item(21, 1).Value = "Biogenic Liquids"
$ws.Cells.Item(21, 2).Value = 2040
$ws.Cells.Item(21, 3).Value = $null
$ws.Cells.Item(21, 4).Value = $null
$ws.Cells.Item(21, 5).Value = $null
$ws.Cells.Item(21, 6).Value = 0.002868182631643872
$ws.Cells.Item(21, 7).Value = 0.00005531636953371903
$ws.Cells.Item(21, 8).Value = 0.00269103731096634
$ws.Cells.Item(21, 9).Value = 0.0025082516179833
$ws.Cells.Item(21, 10).Value = 0.0001125321233566
$ws.Cells.Item(21, 11).Value = 0.0032262798875323

# Row 22: Fossil Liquids (2040)
$ws.Cells.Item(22, 1).Value = "Fossil Liquids"
$ws.Cells.Item(22, 2).Value = 2040
$ws.Cells.Item(22, 3).Value = $null
$ws.Cells.Item(22, 4).Value = $null
$ws.Cells.Item(22, 5).Value = $null
$ws.Cells.Item(22, 6).Value = 0.0182058077380439
$ws.Cells.Item(22, 7).Value = 0.0002569168813649
$ws.Cells.Item(22, 8).Value = 0.0189768366785348
$ws.Cells.Item(22, 9).Value = 0.0107355378940861
$ws.Cells.Item(22, 10).Value = 0.0004975693923245
$ws.Cells.Item(22, 11).Value = 0.0244778202293566

# Row 23: Biomass [Solid] (2040)
$ws.Cells.Item(23, 1).Value = "Biomass [Solid]"
$ws.Cells.Item(23, 2).Value = 2040
$ws.Cells.Item(23, 3).Value = $null
$ws.Cells.Item(23, 4).Value = $null
$ws.Cells.Item(23, 5).Value = 0.001783854964824158
$ws.Cells.Item(23, 6).Value = $null
$ws.Cells.Item(23, 7).Value = $null
$ws.Cells.Item(23, 8).Value = $null
$ws.Cells.Item(23, 9).Value = $null
$ws.Cells.Item(23, 10).Value = $null
$ws.Cells.Item(23, 11).Value = $null

# Row 24: Renewable Energy Carrier (2040)
$ws.Cells.Item(24, 1).Value = "Renewable Energy Carrier"
$ws.Cells.Item(24, 2).Value = 2040
$ws.Cells.Item(24, 3).Value = $null
$ws.Cells.Item(24, 4).Value = $null
$ws.Cells.Item(24, 5).Value = 0.002635382826278307
$ws.Cells.Item(24, 6).Value = $null
$ws.Cells.Item(24, 7).Value = $null
$ws.Cells.Item(24, 8).Value = $null
$ws.Cells.Item(24, 9).Value = $null
$ws.Cells.Item(24, 10).Value = $null
$ws.Cells.Item(24, 11).Value = $null

# Row 25: Overall Demand (2040)
$ws.Cells.Item(25, 1).Value = "Overall Demand"
$ws.Cells.Item(25, 2).Value = 2040
$ws.Cells.Item(25, 3).Value = $null
$ws.Cells.Item(25, 4).Value = 0.004453516139235275
$ws.Cells.Item(25, 5).Value = 0.007822006519660729
$ws.Cells.Item(25, 6).Value = 0.02375433030555169
$ws.Cells.Item(25, 7).Value = 0.000312233250898619
$ws.Cells.Item(25, 8).Value = 0.02166794534336055
$ws.Cells.Item(25, 9).Value = 0.01351885718843164
$ws.Cells.Item(25, 10).Value = 0.0006101015156811001
$ws.Cells.Item(25, 11).Value = 0.0277041001168889

# Row 26: Hydrogen (2050)
$ws.Cells.Item(26, 1).Value = "Hydrogen"
$ws.Cells.Item(26, 2).Value = 2050
$ws.Cells.Item(26, 3).Value = $null
$ws.Cells.Item(26, 4).Value = $null
$ws.Cells.Item(26, 5).Value = $null
$ws.Cells.Item(26, 6).Value = 0.0024695872882565
$ws.Cells.Item(26, 7).Value = $null
$ws.Cells.Item(26, 8).Value = 0.0000001209390442586334
$ws.Cells.Item(26, 9).Value = 0.0002491497128281486
$ws.Cells.Item(26, 10).Value = $null
$ws.Cells.Item(26, 11).Value = $null

# Row 27: Methanol (2050)
$ws.Cells.Item(27, 1).Value = "Methanol"
$ws.Cells.Item(27, 2).Value = 2050
$ws.Cells.Item(27, 3).Value = $null
$ws.Cells.Item(27, 4).Value = $null
$ws.Cells.Item(27, 5).Value = $null
$ws.Cells.Item(27, 6).Value = $null
$ws.Cells.Item(27, 7).Value = $null
$ws.Cells.Item(27, 8).Value = $null
$ws.Cells.Item(27, 9).Value = $null
$ws.Cells.Item(27, 10).Value = $null
$ws.Cells.Item(27, 11).Value = $null

# Row 28: Ammonia (2050)
$ws.Cells.Item(28, 1).Value = "Ammonia"
$ws.Cells.Item(28, 2).Value = 2050
$ws.Cells.Item(28, 3).Value = $null
$ws.Cells.Item(28, 4).Value = 0.004391844017420131
$ws.Cells.Item(28, 5).Value = $null
$ws.Cells.Item(28, 6).Value = $null
$ws.Cells.Item(28, 7).Value = $null
$ws.Cells.Item(28, 8).Value = $null
$ws.Cells.Item(28, 9).Value = $null
$ws.Cells.Item(28, 10).Value = $null
$ws.Cells.Item(28, 11).Value = $null

# Row 29: Synthetic Gases (2050)
$ws.Cells.Item(29, 1).Value = "Synthetic Gases"
$ws.Cells.Item(29, 2).Value = 2050
$ws.Cells.Item(29, 3).Value = $null
$ws.Cells.Item(29, 4).Value = $null
$ws.Cells.Item(29, 5).Value = $null
$ws.Cells.Item(29, 6).Value = 0.000000005419311215647068
$ws.Cells.Item(29, 7).Value = $null
$ws.Cells.Item(29, 8).Value = $null
$ws.Cells.Item(29, 9).Value = 0.000000001264475918242466
$ws.Cells.Item(29, 10).Value = $null
$ws.Cells.Item(29, 11).Value = $null

# Row 30: Biogenic Gases (2050)
$ws.Cells.Item(30, 1).Value = "Biogenic Gases"
$ws.Cells.Item(30, 2).Value = 2050
$ws.Cells.Item(30, 3).Value = $null
$ws.Cells.Item(30, 4).Value = $null
$ws.Cells.Item(30, 5).Value = 0.008523752276313709
$ws.Cells.Item(30, 6).Value = 0.00002188064744053967
$ws.Cells.Item(30, 7).Value = $null
$ws.Cells.Item(30, 8).Value = $null
$ws.Cells.Item(30, 9).Value = 0.000009630866481644564
$ws.Cells.Item(30, 10).Value = $null
$ws.Cells.Item(30, 11).Value = $null

# Row 31: Fossil Gases (2050)
$ws.Cells.Item(31, 1).Value = "Fossil Gases"
$ws.Cells.Item(31, 2).Value = 2050
$ws.Cells.Item(31, 3).Value = $null
$ws.Cells.Item(31, 4).Value = $null
$ws.Cells.Item(31, 5).Value = $null
$ws.Cells.Item(31, 6).Value = 0.00004448656720562416
$ws.Cells.Item(31, 7).Value = $null
$ws.Cells.Item(31, 8).Value = $null
$ws.Cells.Item(31, 9).Value = 0.00003164362716767037
$ws.Cells.Item(31, 10).Value = $null
$ws.Cells.Item(31, 11).Value = $null

# Row 32: Synthetic Liquids (2050)
$ws.Cells.Item(32, 1).Value = "Synthetic Liquids"
$ws.Cells.Item(32, 2).Value = 2050
$ws.Cells.Item(32, 3).Value = $null
$ws.Cells.Item(32, 4).Value = $null
$ws.Cells.Item(32, 5).Value = $null
$ws.Cells.Item(32, 6).Value = 0.00000000002054838133199776
$ws.Cells.Item(32, 7).Value = 0.000000000002037002907176165
$ws.Cells.Item(32, 8).Value = 0.0000000001030654180407414
$ws.Cells.Item(32, 9).Value = 0.00000000004150217993694974
$ws.Cells.Item(32, 10).Value = 0.0000000000007707128990841736
$ws.Cells.Item(32, 11).Value = 0.0000000002936289411462306

# Row 33: Biogenic Liquids (2050)
$ws.Cells.Item(33, 1).Value = "Biogenic Liquids"
$ws.Cells.Item(33, 2).Value = 2050
$ws.Cells.Item(33, 3).Value = $null
$ws.Cells.Item(33, 4).Value = $null
$ws.Cells.Item(33, 5).Value = $null
$ws.Cells.Item(33, 6).Value = 0.0003081190100342624
$ws.Cells.Item(33, 7).Value = 0.0000988406403828583
$ws.Cells.Item(33, 8).Value = 0.003534054646170003
$ws.Cells.Item(33, 9).Value = 0.0006471367414306
$ws.Cells.Item(33, 10).Value = 0.0001447341306725
$ws.Cells.Item(33, 11).Value = 0.0045478180094235

# Row 34: Fossil Liquids (2050)
$ws.Cells.Item(34, 1).Value = "Fossil Liquids"
$ws.Cells.Item(34, 2).Value = 2050
$ws.Cells.Item(34, 3).Value = $null
$ws.Cells.Item(34, 4).Value = $null
$ws.Cells.Item(34, 5).Value = $null
$ws.Cells.Item(34, 6).Value = 0.0012861803472131
$ws.Cells.Item(34, 7).Value = 0.0002317281464704
$ws.Cells.Item(34, 8).Value = 0.0171768198088114
$ws.Cells.Item(34, 9).Value = 0.0019106303723634
$ws.Cells.Item(34, 10).Value = 0.0004282833010402
$ws.Cells.Item(34, 11).Value = 0.0224298349252638

# Row 35: Biomass [Solid] (2050)
$ws.Cells.Item(35, 1).Value = "Biomass [Solid]"
$ws.Cells.Item(35, 2).Value = 2050
$ws.Cells.Item(35, 3).Value = $null
$ws.Cells.Item(35, 4).Value = $null
$ws.Cells.Item(35, 5).Value = 0.001797963928962914
$ws.Cells.Item(35, 6).Value = $null
$ws.Cells.Item(35, 7).Value = $null
$ws.Cells.Item(35, 8).Value = $null
$ws.Cells.Item(35, 9).Value = $null
$ws.Cells.Item(35, 10).Value = $null
$ws.Cells.Item(35, 11).Value = $null

# Row 36: Renewable Energy Carrier (2050)
$ws.Cells.Item(36, 1).Value = "Renewable Energy Carrier"
$ws.Cells.Item(36, 2).Value = 2050
$ws.Cells.Item(36, 3).Value = $null
$ws.Cells.Item(36, 4).Value = $null
$ws.Cells.Item(36, 5).Value = 0.006428867182493626
$ws.Cells.Item(36, 6).Value = $null
$ws.Cells.Item(36, 7).Value = $null
$ws.Cells.Item(36, 8).Value = $null
$ws.Cells.Item(36, 9).Value = $null
$ws.Cells.Item(36, 10).Value = $null
$ws.Cells.Item(36, 11).Value = $null

# Row 37: Overall Demand (2050)
$ws.Cells.Item(37, 1).Value = "Overall Demand"
$ws.Cells.Item(37, 2).Value = 2050
$ws.Cells.Item(37, 3).Value = $null
$ws.Cells.Item(37, 4).Value = 0.004391844017420131
$ws.Cells.Item(37, 5).Value = 0.01675058338777025
$ws.Cells.Item(37, 6).Value = 0.004130259300009623
$ws.Cells.Item(37, 7).Value = 0.0003305687888902612
$ws.Cells.Item(37, 8).Value = 0.02071099549709108
$ws.Cells.Item(37, 9).Value = 0.002848192626249562
$ws.Cells.Item(37, 10).Value = 0.0005730174324834129
$ws.Cells.Item(37, 11).Value = 0.02697765322831624
